$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 12651
$ws1.Range("F3").Value = 7078
$ws1.Range("F6").Value = 444
$ws1.Range("F10").Value = 994
$ws1.Range("F11").Value = 139
$ws1.Range("F13").Value = 996
$ws1.Range("F17").Value = 238
$ws1.Range("F18").Value = 364
$ws1.Range("F19").Value = 20
$ws1.Range("F20").Value = 272
$ws1.Range("F21").Value = 301
$ws1.Range("F22").Value = 46
$ws1.Range("F23").Value = 135
$ws1.Range("F24").Value = 362
$ws1.Range("F25").Value = 5201
$ws1.Range("F26").Value = 67
$ws1.Range("F27").Value = 1412
$ws1.Range("F28").Value = 303
$ws1.Range("F29").Value = 1275
$ws1.Range("F30").Value = 1275
$ws1.Range("F31").Value = 38
$ws1.Range("F32").Value = 11
$ws1.Range("F33").Value = 1326
$ws1.Range("F35").Value = 3
$ws1.Range("F36").Value = 587
$ws1.Range("G31").Value = 198

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 3736
$ws2.Range("F5").Value = 3736
$ws2.Range("F8").Value = 42
$ws2.Range("F20").Value = 47

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 9250
$ws3.Range("F3").Value = 556
$ws3.Range("F4").Value = 1978

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 9250
$ws4.Range("F3").Value = 556
$ws4.Range("F4").Value = 1978
$ws4.Range("F6").Value = 12651
$ws4.Range("F7").Value = 7078
$ws4.Range("F9").Value = 3736
$ws4.Range("F10").Value = 444
$ws4.Range("F12").Value = 995
$ws4.Range("F13").Value = 139
$ws4.Range("F15").Value = 996
$ws4.Range("F19").Value = 238
$ws4.Range("F20").Value = 364
$ws4.Range("F21").Value = 20
$ws4.Range("F22").Value = 272
$ws4.Range("F23").Value = 301
$ws4.Range("F24").Value = 46
$ws4.Range("F29").Value = 362
$ws4.Range("F30").Value = 5201
$ws4.Range("F31").Value = 67
$ws4.Range("F32").Value = 1412
$ws4.Range("F35").Value = 303
$ws4.Range("F37").Value = 1275
$ws4.Range("F38").Value = 1275
$ws4.Range("F39").Value = 1326
$ws4.Range("F41").Value = 587
$ws4.Range("F50").Value = 47

